$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0008834231581628185
$ws.Range("E2").Value = 0.0008834231581628185

# Row 3
$ws.Range("D3").Value = 0.00486079416036151
$ws.Range("E3").Value = 0.00486079416036151

# Row 4
$ws.Range("D4").Value = [double]"3.78699826921037E-15"
$ws.Range("E4").Value = [double]"3.78699826921037E-15"

# Row 5
$ws.Range("D5").Value = [double]"3.789178724624225E-15"
$ws.Range("E5").Value = [double]"3.789178724624225E-15"

# Row 6
$ws.Range("D6").Value = [double]"1.543168793449495E-09"
$ws.Range("E6").Value = [double]"1.543168793449495E-09"

# Row 7
$ws.Range("C7").Value = $true
$ws.Range("D7").Value = 0.8082699960279633
$ws.Range("E7").Value = 0.1917300039720367

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = 0.1295353781629063
$ws.Range("E8").Value = 0.8704646218370937

# Row 9
$ws.Range("D9").Value = 0.9999999999761455
$ws.Range("E9").Value = [double]"2.385447395170104E-11"

# Row 10
$ws.Range("D10").Value = 0.0002063266209142997
$ws.Range("E10").Value = 0.9997936733790856

# Row 11
$ws.Range("D11").Value = 0.9999893352282481
$ws.Range("E11").Value = [double]"1.06647717519337E-05"
$ws.Range("F11").Value = 1.07484757900238
